# Applies the Tue Jun 27 14:55:28 UTC 2023 "Updated cryptos list" data refresh:
# refreshed Price/Volume(1h) figures, and the Dai / ShibaInu / WrappedliquidstakedEther2.0
# and Aave / Algorand rows swapped position in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store numeric-looking values as plain text
# in the source data (e.g. "1.000", "239.46"); force text format first so Excel
# does not silently convert them into numbers when the Value is assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.793.17"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.884.14"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "239.46"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "0.4821"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "0.06547"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "1.996.26"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").Value = "0.07498"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "16.62"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "5.107"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "88.77"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "0.6676"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "30.739.56"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "13.35"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.239.73"
$ws.Range("E18").Value = "  +4.51%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000007630"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "232.74"
$ws.Range("E21").Value = "  +5.71%  "
$ws.Range("D22").Value = "5.303"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "6.186"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "9.315"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "167.76"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").Value = "18.75"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "1.946"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "0.09811"
$ws.Range("E30").Value = "  +6.84%  "
$ws.Range("D31").Value = "4.363"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "0.05079"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "1.217"
$ws.Range("E34").Value = "  +5.23%  "
$ws.Range("D35").Value = "0.7574"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "0.01874"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "2.653"
$ws.Range("D39").Value = "2.094"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "0.9150"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "106.67"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "0.4296"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").Value = "5.816"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "7.370"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "64.78"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1291"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "1.487"
$ws.Range("E48").Value = "  -6.35%  "
$ws.Range("D49").Value = "8.979"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "33.91"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "0.05665"
$ws.Range("E51").Value = "  -0.86%  "
